# Blackjack strategy workbook: add a "split" strategy sheet, recolor the
# conditional-formatting highlights on all three sheets to a plain
# white/grey scheme, and make the new sheet the active one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "split" worksheet after "soft" (becomes 3rd sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSplit = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsSplit.Name = "split"

# ---------------------------------------------------------------------
# 2. Populate the "split" sheet: header row + 10x10 data grid.
# ---------------------------------------------------------------------
$wsSplit.Range("A1").Value = "Player"
$wsSplit.Range("B1").Value = "Dealer2"
$wsSplit.Range("C1").Value = "Dealer3"
$wsSplit.Range("D1").Value = "Dealer4"
$wsSplit.Range("E1").Value = "Dealer5"
$wsSplit.Range("F1").Value = "Dealer6"
$wsSplit.Range("G1").Value = "Dealer7"
$wsSplit.Range("H1").Value = "Dealer8"
$wsSplit.Range("I1").Value = "Dealer9"
$wsSplit.Range("J1").Value = "Dealer10"
$wsSplit.Range("K1").Value = "Dealer11"

for ($r = 2; $r -le 11; $r++) {
    $wsSplit.Cells.Item($r, 1).Value = $r
}
$wsSplit.Range("B2:K11").Value = "No Split"

# ---------------------------------------------------------------------
# 3. Conditional formatting for "split": two 2-rule blocks, same scheme
#    ("No Split" / "Surrender") used on the left (B:F) and right (G:K)
#    halves of the table.
# ---------------------------------------------------------------------
$white = 16777215
$grey = 14277081

$fcBF1 = $wsSplit.Range("B2:F11").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "No Split", 0)
$fcBF1.Interior.Color = $grey
$fcBF2 = $wsSplit.Range("B2:F11").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Surrender", 0)
$fcBF2.Interior.Color = $white

$fcGK1 = $wsSplit.Range("G2:K11").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "No Split", 0)
$fcGK1.Interior.Color = $grey
$fcGK2 = $wsSplit.Range("G2:K11").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Surrender", 0)
$fcGK2.Interior.Color = $white

# ---------------------------------------------------------------------
# 4. Recolor the existing conditional formatting on "hard" and "soft"
#    from the old accent colors to the same plain white/grey scheme.
#    Rule order (by priority) on both sheets is Split, Double, Stand, Hit.
# ---------------------------------------------------------------------
foreach ($sheetName in @("hard", "soft")) {
    $ws = $wb.Worksheets.Item($sheetName)
    if ($sheetName -eq "hard") {
        $rng = $ws.Range("B2:K19")
    } else {
        $rng = $ws.Range("B2:K9")
    }
    $fcs = $rng.FormatConditions
    $fcSplit = $fcs.Item(1)
    $fcDouble = $fcs.Item(2)
    $fcStand = $fcs.Item(3)
    $fcHit = $fcs.Item(4)

    $fcSplit.Interior.Color = $grey
    $fcDouble.Interior.Color = $white
    $fcStand.Interior.Color = $grey
    $fcHit.Interior.Color = $white
}

# ---------------------------------------------------------------------
# 5. Selections / active sheet: "split" becomes the active/selected tab,
#    with D10 selected; "hard" and "soft" keep their own selections.
# ---------------------------------------------------------------------
$wsSplit.Activate()
$wsSplit.Range("D10").Select()
